# Stark MOP bill importer
#
# The workbook used to hold the raw half-hourly (HH) MOP bill data on the
# one and only sheet ("Sheet1"). The importer now splits the workbook into
# two tabs:
#   - "Summary" : a fresh, blank landing sheet (becomes the first tab).
#   - "HH"      : the original bill data (becomes the second, active tab),
#                 with a newly-imported settlement row appended.
#
# NOTE: worksheet object references returned by this host are positional
# (they track a sheet *index*, not a stable identity), so once sheets are
# inserted/reordered we re-fetch fresh references by position rather than
# reusing older variables.

$wb = $excel.ActiveWorkbook

# The single existing sheet keeps its data but needs to end up second, so
# insert a brand-new blank sheet *before* it -- the existing sheet (with
# all its rows/styles/column widths) is left completely undisturbed.
$orig = $wb.Worksheets.Item(1)
$wb.Worksheets.Add($orig) | Out-Null

# Re-fetch by position now that the insert has happened.
$summary = $wb.Worksheets.Item(1)
$hh = $wb.Worksheets.Item(2)

$summary.Name = "Summary"
$hh.Name = "HH"

# Append the new half-hourly settlement row (row 13) to the HH sheet --
# same MPAN / Comms / Type / Site / PS / Settled as row 12, but a
# different settlement period (columns F/G) pulled in by the importer.
$hh.Range("B13").Value = 1472066139971
$hh.Range("B13").NumberFormat = "0000000000000"

$hh.Range("C13").Value = $hh.Range("C12").Value()
$hh.Range("C13").NumberFormat = "0000000000000"

$hh.Range("D13").Value = $hh.Range("D12").Value()

$hh.Range("E13").Value = $hh.Range("E12").Value()
$hh.Range("E13").NumberFormat = "DD/MM/YYYY"

$hh.Range("F13").Value = 29738
$hh.Range("F13").NumberFormat = "DD/MM/YYYY"

$hh.Range("G13").Value = 29767
$hh.Range("G13").NumberFormat = "DD/MM/YYYY"

$hh.Range("H13").Value = 704.3
$hh.Range("I13").Value = 135.81

# Leave the cursor where the importer dropped it and make HH the active tab.
$hh.Range("G14").Select() | Out-Null
$hh.Activate() | Out-Null
